$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.360", "19.30") are preserved exactly as text, matching the
# source data which uses non-standard thousand/decimal separators.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply updated price/volume data scraped from coinranking.com
$ws.Range("D2").Value = "55.554.62"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "2.957.54"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "497.93"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "135.43"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "7.10"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "0.360"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "3.485.72"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "25.65"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "0.0000158"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "55.663.96"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.969.69"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "5.94"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "7.88"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "325.15"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "0.487"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "64.01"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "3.100.36"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "0.0₃0879"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "6.29"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("D30").Value = "6.87"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "19.95"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "1.13"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").Value = "153.16"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "5.63"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "24.82"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("D38").Value = "1.22"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").Value = "0.0650"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").Value = "3.000.62"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "36.54"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "0.646"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "2.141.37"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("D46").Value = "1.33"
$ws.Range("E46").Value = "  -4.53%  "
$ws.Range("D47").Value = "5.76"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("D48").Value = "0.911"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").Value = "0.0232"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "19.30"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "0.0841"
$ws.Range("E51").Value = "  -4.07%  "
